$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header updates: Volume number and report week date range ---
$ws.Range("A8").Value = "Volume 30   Number  51"
$ws.Range("C9").Value = "Report Covering the Week  12/18/2023  Through  12/24/2023"

# --- Weekly crime statistics data updates (rows 14-30) ---
# Row 14
$ws.Range("N14").Value = -81.818181818181

# Row 15
$ws.Range("N15").Value = -71.052631578947

# Row 16
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 19
$ws.Range("H16").Value = 72.727272727272
$ws.Range("I16").Value = 210
$ws.Range("J16").Value = 187
$ws.Range("K16").Value = 12.299465240641
$ws.Range("L16").Value = -0.943396226415
$ws.Range("M16").Value = -50.588235294117
$ws.Range("N16").Value = -89.771066731612

# Row 17
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 17
$ws.Range("E17").Value = -64.705882352941
$ws.Range("F17").Value = 26
$ws.Range("G17").Value = 35
$ws.Range("H17").Value = -25.714285714285
$ws.Range("I17").Value = 403
$ws.Range("J17").Value = 447
$ws.Range("K17").Value = -9.843400447427
$ws.Range("L17").Value = 3.069053708439
$ws.Range("M17").Value = -3.588516746411
$ws.Range("N17").Value = -63.856502242152

# Row 18
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = -33.333333333333
$ws.Range("F18").Value = 13
$ws.Range("G18").Value = 21
$ws.Range("H18").Value = -38.095238095238
$ws.Range("I18").Value = 224
$ws.Range("J18").Value = 244
$ws.Range("K18").Value = -8.196721311475
$ws.Range("L18").Value = -6.276150627615
$ws.Range("M18").Value = -45.098039215686
$ws.Range("N18").Value = -81.993569131832

# Row 19
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = 50
$ws.Range("F19").Value = 30
$ws.Range("G19").Value = 27
$ws.Range("H19").Value = 11.111111111111
$ws.Range("I19").Value = 429
$ws.Range("J19").Value = 409
$ws.Range("K19").Value = 4.889975550122
$ws.Range("L19").Value = 9.160305343511
$ws.Range("M19").Value = 5.925925925925
$ws.Range("N19").Value = -49.883177570093

# Row 20
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 19
$ws.Range("G20").Value = 15
$ws.Range("H20").Value = 26.666666666666
$ws.Range("I20").Value = 134
$ws.Range("J20").Value = 100
$ws.Range("K20").Value = 34
$ws.Range("L20").Value = 47.252747252747
$ws.Range("M20").Value = 16.521739130434
$ws.Range("N20").Value = -75.457875457875

# Row 21
$ws.Range("C21").Value = 26
$ws.Range("D21").Value = 36
$ws.Range("E21").Value = -27.777777777777
$ws.Range("F21").Value = 108
$ws.Range("G21").Value = 110
$ws.Range("H21").Value = -1.818181818181
$ws.Range("I21").Value = 1432
$ws.Range("J21").Value = 1420
$ws.Range("K21").Value = 0.845070422535
$ws.Range("L21").Value = 6.152705707931
$ws.Range("M21").Value = -20.971302428256
$ws.Range("N21").Value = -75.912531539108

# Row 22
$ws.Range("F22").Value = 2
$ws.Range("H22").Value = 100
$ws.Range("I22").Value = 20
$ws.Range("K22").Value = -20
$ws.Range("L22").Value = -25.925925925925
$ws.Range("M22").Value = -41.176470588235

# Row 23
$ws.Range("C23").Value = 5
$ws.Range("D23").Value = 3
$ws.Range("D23").NumberFormat = "#,##0"
$ws.Range("E23").Value = 66.666666666666
$ws.Range("E23").NumberFormat = "#,##0.0;" + [char]34 + "-" + [char]34 + "#,##0.0"
$ws.Range("F23").Value = 18
$ws.Range("G23").Value = 10
$ws.Range("H23").Value = 80
$ws.Range("I23").Value = 241
$ws.Range("J23").Value = 225
$ws.Range("K23").Value = 7.111111111111
$ws.Range("L23").Value = 4.329004329004
$ws.Range("M23").Value = 12.616822429906

# Row 24
$ws.Range("C24").Value = 19
$ws.Range("D24").Value = 22
$ws.Range("E24").Value = -13.636363636363
$ws.Range("F24").Value = 91
$ws.Range("G24").Value = 123
$ws.Range("H24").Value = -26.016260162601
$ws.Range("I24").Value = 1434
$ws.Range("J24").Value = 1592
$ws.Range("K24").Value = -9.924623115577
$ws.Range("L24").Value = 12.118842845973
$ws.Range("M24").Value = 40.313111545988

# Row 25
$ws.Range("C25").Value = 15
$ws.Range("D25").Value = 11
$ws.Range("E25").Value = 36.363636363636
$ws.Range("F25").Value = 48
$ws.Range("G25").Value = 36
$ws.Range("H25").Value = 33.333333333333
$ws.Range("I25").Value = 634
$ws.Range("J25").Value = 646
$ws.Range("K25").Value = -1.857585139318
$ws.Range("L25").Value = 24.803149606299
$ws.Range("M25").Value = -33.403361344537

# Row 26
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = "0"

# Row 27
$ws.Range("C27").Value = 3
$ws.Range("C27").NumberFormat = "#,##0"
$ws.Range("E27").Value = 200
$ws.Range("F27").Value = 4
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 100
$ws.Range("I27").Value = 47
$ws.Range("J27").Value = 43
$ws.Range("K27").Value = 9.302325581395
$ws.Range("L27").Value = -41.25

# Row 28
$ws.Range("D28").Value = 1
$ws.Range("D28").NumberFormat = "#,##0"
$ws.Range("E28").Value = -100
$ws.Range("E28").NumberFormat = "#,##0.0;" + [char]34 + "-" + [char]34 + "#,##0.0"
$ws.Range("F28").Value = 1
$ws.Range("G28").Value = 3
$ws.Range("H28").Value = -66.666666666666
$ws.Range("J28").Value = 34
$ws.Range("K28").Value = 0
$ws.Range("N28").Value = -87.5

# Row 29
$ws.Range("D29").Value = 1
$ws.Range("D29").NumberFormat = "#,##0"
$ws.Range("E29").Value = -100
$ws.Range("E29").NumberFormat = "#,##0.0;" + [char]34 + "-" + [char]34 + "#,##0.0"
$ws.Range("F29").Value = 1
$ws.Range("G29").Value = 3
$ws.Range("H29").Value = -66.666666666666
$ws.Range("J29").Value = 27
$ws.Range("K29").Value = -7.407407407407
$ws.Range("N29").Value = -89.583333333333

# Row 30
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "0"
